# Applies the text/formatting fixes described in the commit diff:
#  1. Slide 3 - "Built specifically for Windows 10 laptops..." -> drop "10"
#  2. Slide 6 - "Designed for Windows 10 laptops..." -> drop "10"
#  3. Slide 7 - Content placeholder shape nudged 1 EMU to the left (525717 -> 525716)
#  4. Slide 8 - "Focus only on Linux platforms (e.g., Aircrack-ng)." -> trimmed, merged to one run
#  5. Slide 8 - Table-like text lines lose their leading/inline "|" pipe characters

$p = $ppt.ActivePresentation

# --- Change 1: Slide 3, Content Placeholder 2, paragraph 2 ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$para3_2 = $tr3.Paragraphs(2, 1)
$run3_2 = $para3_2.Runs(1, 1)
$run3_2.Text = " - Built specifically for Windows laptops to ensure wide usability.  "

# --- Change 2: Slide 6, Content Placeholder 2, paragraph 10 ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$para6_10 = $tr6.Paragraphs(10, 1)
$run6_10 = $para6_10.Runs(1, 1)
$run6_10.Text = "   - Designed for Windows laptops—used by millions globally. "

# --- Change 3: Slide 7, Content Placeholder 2, nudge x offset by 1 EMU ---
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$sh7.Left = 525716 / 12700.0

# --- Change 4: Slide 8, Content Placeholder 2, paragraph 3 (merge 3 runs into 1) ---
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
$tr8 = $sh8.TextFrame.TextRange

$para8_3 = $tr8.Paragraphs(3, 1)
$all8_3 = $para8_3.Characters(1, $para8_3.Length)
$all8_3.Text = "- Focus only on Linux platforms.  "

# --- Change 5: Slide 8, Content Placeholder 2, paragraphs 7-10 (strip "|" table formatting) ---
$para8_7 = $tr8.Paragraphs(7, 1)
$all8_7 = $para8_7.Characters(1, $para8_7.Length)
$all8_7.Text = "Feature                         Our Tool         Existing Tools "

$para8_8 = $tr8.Paragraphs(8, 1)
$all8_8 = $para8_8.Characters(1, $para8_8.Length)
$all8_8.Text = "Platform Support        Windows       Mostly Linux     "

$para8_9 = $tr8.Paragraphs(9, 1)
$all8_9 = $para8_9.Characters(1, $para8_9.Length)
$all8_9.Text = "Automatic Blocking    Yes                  No/Manual        "

$para8_10 = $tr8.Paragraphs(10, 1)
$all8_10 = $para8_10.Characters(1, $para8_10.Length)
$all8_10.Text = "GUI                                 Yes                   CLI-Based        "
